# Added transportation support for other countries
# Appends Estonia / Latvia / Lithuania / Finland shipping rows (with a
# boolean "active" flag in column A) plus a "Country:" label row to the
# "translation" sheet, which drives the site's i18n lookup table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translation")

# --- Row 104: Estonia -------------------------------------------------
$ws.Range("A104").Value = $true
$ws.Range("B104").Value = "Estonia"
$ws.Range("C104").Value = "Эстония"

# --- Row 105: Latvia ----------------------------------------------------
$ws.Range("A105").Value = $true
$ws.Range("C105").Value = "Латвия"

# --- Row 106: Lithuania ---------------------------------------------------
$ws.Range("A106").Value = $true
$ws.Range("C106").Value = "Литва"

# --- Row 107: Finland (ENG typo "Finalnd" kept verbatim from source) ------
$ws.Range("A107").Value = $true
$ws.Range("C107").Value = "Финляндия"

# --- English names for Latvia / Lithuania / Finland -----------------------
$ws.Range("B105").Value = "Latvia"
$ws.Range("B106").Value = "Lithuania"
$ws.Range("B107").Value = "Finalnd"

# --- Estonian names for Estonia / Latvia / Lithuania / Finland ------------
$ws.Range("D104").Value = "Eesti"
$ws.Range("D105").Value = "Läti"
$ws.Range("D106").Value = "Leedu"
$ws.Range("D107").Value = "Soome"

# --- Row 108: "Country:" label (no column-A flag on this row) -------------
$ws.Range("B108").Value = "Country:"
$ws.Range("D108").Value = "Riik:"
$ws.Range("C108").Value = "Страна:"

# Match the author's final selection/viewport on save.
$ws.Range("C108").Select()
